$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New test case TC-118 appended as rows 31 (header) and 32 (data), mirroring
# the existing TC-119 block (rows 25/26), plus two new trailing columns:
#   Z  -> "Hotel select title" / "Adactin.com - Select Hotel" (copy of O28/O29)
#   AA -> "Results found" header only (new column, no data yet)
# ---------------------------------------------------------------------------

# ---- Row 31: header row (yellow-filled, like row 25/28 etc.) -------------
$headerVals = @{
    "A31" = "TC-118"
    "B31" = "username"
    "C31" = "password"
    "D31" = "location"
    "E31" = "hotel"
    "F31" = "Check In Date"
    "G31" = "Room Type"
    "H31" = "No of rooms"
    "I31" = "Check out Date"
    "J31" = "Adults per Room"
    "K31" = "Children per Room  "
    "L31" = "Login page title"
    "M31" = "Hotel search title"
    "N31" = "Select page title"
    "O31" = "Booking page title"
    "P31" = "Logout page title"
    "Q31" = "message"
    "R31" = "First Name"
    "S31" = "LastName"
    "T31" = "Address"
    "U31" = "Credit num"
    "V31" = "Card type"
    "W31" = "Expiry month"
    "X31" = "Expiry year"
    "Y31" = "CVV number"
    "Z31" = "Hotel select title"
    "AA31" = "Results found"
}
foreach ($addr in $headerVals.Keys) {
    $ws.Range($addr).Value = $headerVals[$addr]
    $ws.Range($addr).Interior.Color = 65535
}

# ---- Row 32: data row -----------------------------------------------------
$dataVals = @{
    "A32" = "TC-118"
    "B32" = "reyaz0806"
    "C32" = "reyaz123"
    "D32" = "Sydney"
    "E32" = "Hotel Creek"
    "F32" = "19/01/2025"
    "G32" = "Standard"
    "H32" = "1 - One"
    "I32" = "20/01/2025"
    "J32" = "1 - One"
    "K32" = "0 - None"
    "L32" = "Adactin.com - Hotel Reservation System"
    "M32" = "Adactin.com - Search Hotel"
    "N32" = "Adactin.com - Select Hotel"
    "O32" = "Adactin.com - Book A Hotel"
    "P32" = "Adactin.com - Logout"
    "Q32" = "Cancel Selected"
    "R32" = "Test"
    "S32" = "Data"
    "T32" = "Hyderabad"
    "U32" = "1234567812345678"
    "V32" = "Master Card"
    "W32" = "March"
    "X32" = "2026"
    "Y32" = "000"
    "Z32" = "Adactin.com - Select Hotel"
}
foreach ($addr in $dataVals.Keys) {
    $ws.Range($addr).Value = $dataVals[$addr]
}

# quotePrefix'd date-like text cells (match F26/I26 pattern)
$ws.Range("F32").NumberFormat = "mm/dd/yyyy"
$ws.Range("F32").Value = "19/01/2025"
$ws.Range("I32").NumberFormat = "mm/dd/yyyy"
$ws.Range("I32").Value = "20/01/2025"

# wrap-text cell (matches Q26 pattern)
$ws.Range("Q32").WrapText = $true

# quotePrefix'd numeric-looking text cells (match U26/X26/Y26 pattern)
$ws.Range("U32").Value = "'1234567812345678"
$ws.Range("X32").Value = "'2026"
$ws.Range("Y32").Value = "'000"

$ws.Rows.Item(32).RowHeight = 28.8

# ---- Column widths for the two new trailing columns (Y=25, Z=26) ---------
$ws.Columns.Item(25).ColumnWidth = 12.109375
$ws.Columns.Item(26).ColumnWidth = 17.6640625

# ---- Sheet view: scrolled position + active selection ---------------------
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 18
$ws.Range("AC36").Select()
